$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.424.50"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.842.19"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.35"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6277"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07406"
$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07721"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "1.846.18"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.964"
$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6706"
$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001042"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.73"
$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.252"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "29.461.28"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.05"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("E22").Value = "  -3.26%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.94"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.477"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1348"
$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.32"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07291"
$ws.Range("E28").Value = "  +11.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value = "  +5.08%  "

$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.160"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.816"
$ws.Range("E34").Value = "  -1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7144"
$ws.Range("E35").Value = "  +1.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.584"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01838"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.790"
$ws.Range("E38").Value = "  -1.96%  "

$ws.Range("D39").Value = "1.232.86"
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.803"
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9558"
$ws.Range("E41").Value = "  +2.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "2.001.32"
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.17"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.34"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("E46").Value = "  +2.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.701"
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.964"
$ws.Range("E48").Value = "  -1.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.892"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("E50").Value = "  -2.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3889"
$ws.Range("E51").Value = "  -1.71%  "
